$d = $word.ActiveDocument

$introText = "Introducci" + [char]0x00F3 + "n."

# Work from the bottom of the affected region upward so that each deletion
# cannot invalidate paragraph indices resolved for steps still to come.

# 4) Remove the blank paragraph right after "Introduccion.".
$pIntro = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.TrimEnd("`r") -eq $introText) {
        $pIntro = $i
    }
}
$pAfterIntro = $d.Paragraphs.Item($pIntro + 1)
$pAfterIntro.Range.Delete()

# 3) Remove the "2. Justificacion..." paragraph and the blank paragraph that
#    follows it (directly before the page-break paragraph).
$pJustificacion = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("2. Justificaci")) {
        $pJustificacion = $i
    }
}
$rJustBlock = $d.Range($d.Paragraphs.Item($pJustificacion).Range.Start, $d.Paragraphs.Item($pJustificacion + 1).Range.End)
$rJustBlock.Delete()

# 2) Empty the "1. Analisis..." paragraph but keep the paragraph (and its
#    formatting) in place.
$pAnalisis = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("1. An")) {
        $pAnalisis = $i
    }
}
$pA = $d.Paragraphs.Item($pAnalisis)
$rClear = $d.Range($pA.Range.Start, $pA.Range.End - 1)
$rClear.Text = ""

# 1) Remove the six paragraphs that precede the "1. Analisis..." paragraph:
#    the brief intro line, a blank line, the quoted requirement paragraph,
#    a blank line, the "Ademas..." paragraph and a trailing blank line.
$firstToDelete = $pAnalisis - 6
$lastToDelete = $pAnalisis - 1
$rIntroBlock = $d.Range($d.Paragraphs.Item($firstToDelete).Range.Start, $d.Paragraphs.Item($lastToDelete).Range.End)
$rIntroBlock.Delete()
